$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 22: "PWM bits" label, LOG(B3,2) formula (same style as other
# formula cells such as B17), and "bit" unit label
$ws.Range("A22").Value = "PWM bits"
$ws.Range("B22").Formula = "=LOG(B3,2)"
$ws.Range("C22").Value = "bit"

# Match the formatting used by the other computed cells (e.g. B17)
$ws.Range("B17").Copy()
$ws.Range("B22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the saved selection to B16, matching the authored workbook state
$ws.Range("B16").Select()
